$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three new columns (E, F, G) before the existing "Operandos" column,
# which pushes it from E to H.
$ws.Columns("E:G").Insert()

# --- Table 1 header (row 4): Opcode / Funct / Code ---
$ws.Range("E4").Value = "Opcode"
$ws.Range("F4").Value = "Funct"
$ws.Range("G4").Value = "Code"

# --- Table 1 first data row (row 5, "add") gets sample values ---
$ws.Range("D5").Value = "R1"
$ws.Range("E5").Value = "0x00"
$ws.Range("F5").Value = "0x20"
$ws.Range("G5").Value = "0x60"
$ws.Range("H5").Value = "Reg, Reg, Reg"

# Fix up the selection/active cell (matches the author's final cursor spot)
$ws.Range("J5").Select() | Out-Null
